$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "jMAkQ849"
$ws.Range("B2").Value = 23102507
$ws.Range("C2").Value = "qdtflzf22"
$ws.Range("D2").Value = "x!Sn5P7#"
$ws.Range("F2").Value = "BsWeWEpS"
$ws.Range("G2").Value = "WqxK"
